# Applies the update described by the diff:
#  - Inserts 3 new product rows (CETAL, OTRIVIN, PANADOL) into the table
#  - Updates the running total and shifts the summary/footer rows down
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1) Insert a new row for "CETAL 250MG/5ML 60ML SUSP" above row 5
#    (which currently holds INJECTMOL), copying the row formatting
#    of the existing data rows so styles/merges match.
# ---------------------------------------------------------------
$ws.Rows(5).Insert()
$ws.Range("A6:N6").Copy()
$ws.Range("A5:N5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B5:G5").Merge()
$ws.Range("H5:K5").Merge()
$ws.Range("L5:M5").Merge()
$ws.Rows(5).RowHeight = 25.5

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "CETAL 250MG/5ML 60ML SUSP"
$ws.Range("H5").Value = "17:0"
$ws.Range("L5").Value = 31
$ws.Range("N5").Value = "1:0"

# ---------------------------------------------------------------
# 2) Insert two new rows above the row holding "جهاز محلول"
#    (currently row 9, since CETAL pushed it down by one) for
#    "OTRIVIN ..." and "PANADOL ..."
# ---------------------------------------------------------------
$ws.Rows(9).Insert()
$ws.Range("A8:N8").Copy()
$ws.Range("A9:N9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Rows(9).RowHeight = 24.75

$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "OTRIVIN 0.05% PEDIATRIC NASAL DROPS 15 ML"
$ws.Range("H9").Value = "5:0"
$ws.Range("L9").Value = 24
$ws.Range("N9").Value = "1:0"

$ws.Rows(10).Insert()
$ws.Range("A9:N9").Copy()
$ws.Range("A10:N10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Rows(10).RowHeight = 25.5

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "PANADOL ADVANCE 500 MG 48 TABLETS"
$ws.Range("H10").Value = "3:3"
$ws.Range("L10").Value = 23
$ws.Range("N10").Value = "0:0"

# ---------------------------------------------------------------
# 3) Renumber the remaining (pre-existing) product rows, which have
#    all shifted down (INJECTMOL, LEVANIC, ORS by one; جهاز محلول,
#    حفاضات, سرنجات, كالونا by three)
# ---------------------------------------------------------------
$ws.Range("A6").Value = 3
$ws.Range("A7").Value = 4
$ws.Range("A8").Value = 5
$ws.Range("A11").Value = 8
$ws.Range("A12").Value = 9
$ws.Range("A13").Value = 10
$ws.Range("A14").Value = 11

# ---------------------------------------------------------------
# 4) Update the running-total cell (now on row 15) with the new sum
# ---------------------------------------------------------------
$ws.Range("K15").Value = 301.67

# ---------------------------------------------------------------
# 5) Make sure every data row keeps its original (auto-fit) height
# ---------------------------------------------------------------
$ws.Rows(4).RowHeight = 24.75
$ws.Rows(5).RowHeight = 25.5
$ws.Rows(6).RowHeight = 24.75
$ws.Rows(7).RowHeight = 25.5
$ws.Rows(8).RowHeight = 25.5
$ws.Rows(9).RowHeight = 24.75
$ws.Rows(10).RowHeight = 25.5
$ws.Rows(11).RowHeight = 24.75
$ws.Rows(12).RowHeight = 25.5
$ws.Rows(13).RowHeight = 25.5
$ws.Rows(14).RowHeight = 24.75
